# This script regenerates the "handback-status" report with new file
# identifiers / timestamps, mirroring what the CI tooling does each run:
#   - the en-US source file guid changes from 8ae8fb6d-... to 03824bdb-...
#   - the en-US source file guid changes from cf23c1d6-... to ffffb3aabc47-...
#   - the handoff content hash changes (shared for both zh-cn and de-de xlf)
#   - the handoff/handback timestamps advance

$wb = $excel.ActiveWorkbook

$oldUuid1 = "8ae8fb6d-2ecf-41b9-9f1d-bb69b4d2f4bc"
$oldUuid2 = "cf23c1d6-efca-4b4a-adbe-bf046a436cb0"
$newUuid1 = "03824bdb-d6e6-4451-b5b5-471410216c3e"
$newUuid2 = "ffffb3aabc47-f91b-4db8-9836-c79007cce269"

$newUuid1Md = "$newUuid1.md"
$newUuid2Md = "$newUuid2.md"

$newHash = "2239120a0a49db1a9e1ad4d040b293a8da2e7ba3"

$newZhCnXlf = "$newUuid1.$newHash.zh-cn.xlf"
$newDeDeXlf = "$newUuid1.$newHash.de-de.xlf"

$zhCnHandoffTime = "2016-03-12 02:43:49"
$zhCnHandbackTime = "2016-03-12 02:44:05"
$deDeHandoffTime = "2016-03-12 02:43:52"
$deDeHandbackTime = "2016-03-12 02:44:10"

# ---------------------------------------------------------------------
# Sheet "Overview": update the two source-file cells and their hyperlinks
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newUuid1Md
$wsOverview.Range("A3").Value = $newUuid2Md

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address(0, 0)
    if ($addr -eq "A2") {
        $h.TextToDisplay = $newUuid1Md
    } elseif ($addr -eq "A3") {
        $h.TextToDisplay = $newUuid2Md
    }
}

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de" share the same layout:
#   A col = source file, B col = file extension, D/G cols = handoff file
#   (xlf), E col = handoff datetime, F col = source file (repeat),
#   H col = handback datetime
# Row 2 corresponds to uuid1, row 3 to uuid2; both rows now point at the
# same (new) handoff file name / timestamps.
# ---------------------------------------------------------------------
function Update-LangSheet($ws, $xlfName, $handoffTime, $handbackTime) {
    $ws.Range("A2").Value = $newUuid1Md
    $ws.Range("F2").Value = $newUuid1Md
    $ws.Range("D2").Value = $xlfName
    $ws.Range("G2").Value = $xlfName
    $ws.Range("E2").Value = $handoffTime
    $ws.Range("H2").Value = $handbackTime

    $ws.Range("A3").Value = $newUuid2Md
    $ws.Range("F3").Value = $newUuid2Md
    $ws.Range("D3").Value = $xlfName
    $ws.Range("G3").Value = $xlfName
    $ws.Range("E3").Value = $handoffTime
    $ws.Range("H3").Value = $handbackTime

    $map = @{
        "A2" = $newUuid1Md;
        "F2" = $newUuid1Md;
        "D2" = $xlfName;
        "G2" = $xlfName;
        "A3" = $newUuid2Md;
        "F3" = $newUuid2Md;
        "D3" = $xlfName;
        "G3" = $xlfName;
    }

    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address(0, 0)
        if ($map.ContainsKey($addr)) {
            $h.TextToDisplay = $map[$addr]
        }
    }
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $wsZhCn $newZhCnXlf $zhCnHandoffTime $zhCnHandbackTime

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LangSheet $wsDeDe $newDeDeXlf $deDeHandoffTime $deDeHandbackTime
